# Daily attendance processing - clear "Recorded By" values in column G
# and shrink column G width now that it's empty.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,22,23,24,25,26,27,28,29,32,33,34,36,37,38,43,44,45,46,47,48,49,50,53,54,55,57,58,59,64,65,66,67,68,69,70,71,74,75,76,78,79,80,85,86,87,88,89,90,91,92,93,94,95,96,97,98,99,100,105,106,107,108,109,110,111,112,113,114,115,116,117,118,119,120,125,126,127,128,129,130,131,132,133,134,135,136,137,138,139,140,145,146,147,148,149,150,151,152,153,154,155,156,157,158,159,160,165,166,167,168,169,170,171,172,173,174,175,176,177,178,179,180,185,186,187,188,189,190,191,192,195,196,197,199,200,201,206,207,208,209,210,211,212,213,216,217,218,220,221,222,227,228,229,230,231,232,233,234,237,238,239,241,242,243)

foreach ($r in $rows) {
    $ws.Cells.Item($r, 7).Value = ""
}

$ws.Columns.Item(7).ColumnWidth = 12.1666666666667
